$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - "Save" with same formatting as G1 (bold, bordered, centered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data cells H2:H4
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
